$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J, matching the style of existing headers (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-26
$values = @(
    @(9, 9),
    @(11, 11),
    @(4, 4),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(7, 8),
    @(6, 6),
    @(5, 5),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
